$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.892.62'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '2.747.57'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = "'573.91"
$ws.Range('E5').Value = '  -1.26%  '
$ws.Range('D6').Value = "'157.13"
$ws.Range('E6').Value = '  +1.22%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -1.37%  '
$ws.Range('E9').Value = '  -3.05%  '
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('E11').Value = '  -2.10%  '
$ws.Range('D12').Value = "'5.46"
$ws.Range('E12').Value = '  -18.87%  '
$ws.Range('D13').Value = '3.231.74'
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('D14').Value = "'26.46"
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('D15').Value = '63.726.46'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('E16').Value = '  -2.45%  '
$ws.Range('D17').Value = '2.750.72'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('E18').Value = '  +0.71%  '
$ws.Range('D19').Value = "'4.78"
$ws.Range('E19').Value = '  -2.18%  '
$ws.Range('D20').Value = "'354.29"
$ws.Range('E20').Value = '  -1.99%  '
$ws.Range('E21').Value = '  -3.44%  '
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').Value = "'65.20"
$ws.Range('E24').Value = '  -2.06%  '
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('D27').Value = "'8.39"
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('E28').Value = '  -2.39%  '
$ws.Range('E29').Value = '  -4.11%  '
$ws.Range('D30').Value = "'6.94"
$ws.Range('E30').Value = '  -2.32%  '
$ws.Range('D31').Value = "'169.19"
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('D32').Value = "'1.21"
$ws.Range('E32').Value = '  -6.20%  '
$ws.Range('E33').Value = '  -2.04%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('E36').Value = '  -0.17%  '
$ws.Range('D37').Value = "'1.77"
$ws.Range('E37').Value = '  -2.38%  '
$ws.Range('E38').Value = '  -3.80%  '
$ws.Range('E39').Value = '  +6.05%  '
$ws.Range('E40').Value = '  -3.11%  '
$ws.Range('D41').Value = "'324.31"
$ws.Range('E41').Value = '  -6.77%  '
$ws.Range('D42').Value = "'38.89"
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('D43').Value = "'21.19"
$ws.Range('E43').Value = '  -3.27%  '
$ws.Range('D44').Value = "'0.0586"
$ws.Range('E44').Value = '  -0.71%  '
$ws.Range('D45').Value = "'21.22"
$ws.Range('E45').Value = '  -2.79%  '
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('D47').Value = "'134.50"
$ws.Range('E47').Value = '  -2.29%  '
$ws.Range('E48').Value = '  -4.22%  '
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').Value = "'11.04"
$ws.Range('E51').Value = '  +0.44%  '
